$wb = $excel.ActiveWorkbook
$croatia = $wb.Worksheets.Item("Croatia")

# Create the new Greece sheet as a copy of Croatia, placed right after it
$croatia.Copy([System.Reflection.Missing]::Value, $croatia)
$greece = $wb.Worksheets.Item($croatia.Index + 1)
$greece.Name = "Greece"

# Update the values (set B4 first so the shared-string table order matches)
$greece.Range("B4").Value = "NGC-4119/T3167/T3166"
$greece.Range("B2").Value = "Greece Market"

# Adjust selections: Croatia's old selection becomes a "select all" with no active sheet tab
$croatia.Select() | Out-Null
$croatia.Range("A1:XFD1048576").Select() | Out-Null

# Greece becomes the active/selected tab with B2 as the active cell
$greece.Select() | Out-Null
$greece.Range("B2").Select() | Out-Null
